$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade the last few homework 1 entries that were still blank:
#   row 8  - Ismail, Noha    -> good
#   row 9  - Laderman, Eric  -> excellent
#   row 16 - Yong, Luok Wen  -> excellent
# These cells currently carry the "ungraded" (red) fill; clear that
# formatting the same way it was already cleared on B4 by copying its
# format (border kept, fill cleared) before writing the new values.
$ws.Range("B4").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B8").Value = "good"
$ws.Range("B9").Value = "excellent"
$ws.Range("B16").Value = "excellent"

# Move the active selection to B17, matching where the cursor was left.
[void]$ws.Range("B17").Select()
